$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Pre-format E2:G39 as text so date-like / number-like strings
# ("10/07/2021", "5.796,00", etc.) are stored as literal text instead
# of being auto-converted to Excel date/number serials.
$ws.Range("E2:G39").NumberFormat = "@"

$ws.Range("E2").Value = "R120637"
$ws.Range("F2").Value = "5.796,00"
$ws.Range("G2").Value = "10/07/2021"
$ws.Range("I2").Value = 85.724999999999994
$ws.Range("E3").Value = "R120639"
$ws.Range("F3").Value = "5.796,00"
$ws.Range("G3").Value = "10/07/2021"
$ws.Range("I3").Value = 82.224000000000004
$ws.Range("E4").Value = "R120637"
$ws.Range("F4").Value = "5.796,00"
$ws.Range("G4").Value = "25/07/2021"
$ws.Range("I4").Value = 90.701999999999998
$ws.Range("E5").Value = "R120639"
$ws.Range("F5").Value = "5.796,00"
$ws.Range("G5").Value = "25/07/2021"
$ws.Range("I5").Value = 80.936000000000007
$ws.Range("E6").Value = "R120637"
$ws.Range("I6").Value = 83.784999999999997
$ws.Range("E7").Value = "R120639"
$ws.Range("I7").Value = 83.674000000000007
$ws.Range("E8").Value = "R120637"
$ws.Range("F8").Value = "9.660,00"
$ws.Range("G8").Value = "10/10/2021"
$ws.Range("I8").Value = 85.724999999999994
$ws.Range("E9").Value = "R120639"
$ws.Range("F9").Value = "9.660,00"
$ws.Range("G9").Value = "10/10/2021"
$ws.Range("I9").Value = 82.224000000000004
$ws.Range("E10").Value = "R120637"
$ws.Range("F10").Value = "7.728,00"
$ws.Range("G10").Value = "10/11/2021"
$ws.Range("I10").Value = 91.055999999999997
$ws.Range("E11").Value = "R120639"
$ws.Range("F11").Value = "7.728,00"
$ws.Range("G11").Value = "10/11/2021"
$ws.Range("I11").Value = 80.936000000000007
$ws.Range("E12").Value = "R120637"
$ws.Range("F12").Value = "9.660,00"
$ws.Range("G12").Value = "05/12/2021"
$ws.Range("I12").Value = 85.724999999999994
$ws.Range("E13").Value = "R120639"
$ws.Range("F13").Value = "9.660,00"
$ws.Range("G13").Value = "05/12/2021"
$ws.Range("I13").Value = 82.224000000000004
$ws.Range("E14").Value = "R120637"
$ws.Range("F14").Value = "9.660,00"
$ws.Range("G14").Value = "05/01/2022"
$ws.Range("I14").Value = 91.055999999999997
$ws.Range("E15").Value = "R120639"
$ws.Range("F15").Value = "9.660,00"
$ws.Range("G15").Value = "05/01/2022"
$ws.Range("I15").Value = 80.936000000000007
$ws.Range("E16").Value = "R120637"
$ws.Range("F16").Value = "9.660,00"
$ws.Range("G16").Value = "05/02/2022"
$ws.Range("I16").Value = 83.784999999999997
$ws.Range("E17").Value = "R120639"
$ws.Range("F17").Value = "9.660,00"
$ws.Range("G17").Value = "05/02/2022"
$ws.Range("I17").Value = 83.674000000000007
$ws.Range("E18").Value = "R120637"
$ws.Range("F18").Value = "7.728,00"
$ws.Range("G18").Value = "05/03/2022"
$ws.Range("I18").Value = 85.724999999999994
$ws.Range("E19").Value = "R120639"
$ws.Range("F19").Value = "7.728,00"
$ws.Range("G19").Value = "05/03/2022"
$ws.Range("I19").Value = 82.224000000000004
$ws.Range("E20").Value = "R120637"
$ws.Range("F20").Value = "7.728,00"
$ws.Range("G20").Value = "05/04/2022"
$ws.Range("I20").Value = 91.055999999999997
$ws.Range("E21").Value = "R120639"
$ws.Range("F21").Value = "7.728,00"
$ws.Range("G21").Value = "05/04/2022"
$ws.Range("I21").Value = 80.936000000000007
$ws.Range("E22").Value = "R120637"
$ws.Range("F22").Value = "7.728,00"
$ws.Range("G22").Value = "05/05/2022"
$ws.Range("I22").Value = 83.784999999999997
$ws.Range("E23").Value = "R120639"
$ws.Range("F23").Value = "7.728,00"
$ws.Range("G23").Value = "05/05/2022"
$ws.Range("I23").Value = 83.674000000000007
$ws.Range("E24").Value = "R120637"
$ws.Range("F24").Value = "9.660,00"
$ws.Range("G24").Value = "05/12/2021"
$ws.Range("I24").Value = 85.724999999999994
$ws.Range("E25").Value = "R120639"
$ws.Range("F25").Value = "9.660,00"
$ws.Range("G25").Value = "05/12/2021"
$ws.Range("I25").Value = 82.224000000000004
$ws.Range("E26").Value = "R120637"
$ws.Range("G26").Value = "05/01/2022"
$ws.Range("I26").Value = 91.055999999999997
$ws.Range("E27").Value = "R120639"
$ws.Range("G27").Value = "05/01/2022"
$ws.Range("I27").Value = 80.936000000000007
$ws.Range("E28").Value = "R120637"
$ws.Range("G28").Value = "05/02/2022"
$ws.Range("I28").Value = 83.784999999999997
$ws.Range("E29").Value = "R120639"
$ws.Range("G29").Value = "05/02/2022"
$ws.Range("I29").Value = 83.674000000000007
$ws.Range("E30").Value = "R120637"
$ws.Range("F30").Value = "7.728,00"
$ws.Range("G30").Value = "05/03/2022"
$ws.Range("I30").Value = 85.724999999999994
$ws.Range("E31").Value = "R120639"
$ws.Range("F31").Value = "7.728,00"
$ws.Range("G31").Value = "05/03/2022"
$ws.Range("I31").Value = 82.224000000000004
$ws.Range("E32").Value = "R120637"
$ws.Range("F32").Value = "9.660,00"
$ws.Range("G32").Value = "05/12/2021"
$ws.Range("I32").Value = 85.724999999999994
$ws.Range("E33").Value = "R120639"
$ws.Range("F33").Value = "9.660,00"
$ws.Range("G33").Value = "05/12/2021"
$ws.Range("I33").Value = 82.224000000000004
$ws.Range("E34").Value = "R120637"
$ws.Range("G34").Value = "05/01/2022"
$ws.Range("I34").Value = 91.055999999999997
$ws.Range("E35").Value = "R120639"
$ws.Range("G35").Value = "05/01/2022"
$ws.Range("I35").Value = 80.936000000000007
$ws.Range("E36").Value = "R120637"
$ws.Range("G36").Value = "05/02/2022"
$ws.Range("I36").Value = 83.784999999999997
$ws.Range("E37").Value = "R120639"
$ws.Range("G37").Value = "05/02/2022"
$ws.Range("I37").Value = 83.674000000000007
$ws.Range("E38").Value = "R120637"
$ws.Range("F38").Value = "7.728,00"
$ws.Range("G38").Value = "05/03/2022"
$ws.Range("I38").Value = 85.724999999999994
$ws.Range("E39").Value = "R120639"
$ws.Range("F39").Value = "7.728,00"
$ws.Range("G39").Value = "05/03/2022"
$ws.Range("I39").Value = 82.224000000000004

# Restore the default "Normal" style so no stray per-cell formatting
# is left behind (matches the original workbook's styling).
$ws.Range("E2:G39").Style = "Normal"

# Update sheet view state (scroll position / active selection).
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("L46").Select()

Write-Host "Applied Worldclass Industries reference range updates"
